$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value()
    $cell.Value = ($old -replace "_old$", "_FV2404")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $old = $cell.Value()
    $cell.Value = ($old -replace "_new$", "_FV2410")
}

$ws.Activate()
$ws.Range("A2").Select()
$win = $excel.ActiveWindow()
$win.FreezePanes = $true

$los = $ws.ListObjects()
$tbl = $los.Add(1, $ws.Range("A1:U70"), $null, 1)
